# Add new rows of health-tracker data (rows 3-8), including a date field in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing data row (row 2) down onto the new rows,
# so the new cells reuse the same styles (bordered/bold id column, date-formatted
# date column) instead of creating brand-new style entries.
$ws.Range("A2:M2").Copy()
$ws.Range("A3:M8").PasteSpecial(-4122)

$rows = @(
    @{ n=3; A=1; B="no"; C=2; D=3; E=3; F="yes"; G=43711.36360478804; H=3; I=45; J=5; K="no"; L="looking at phone"; M=12 },
    @{ n=4; A=2; B="no"; C=3; D=3; E=3; F="yes"; G=43711.36394776197; H=3; I=45; J=5; K="no"; L="looking at phone"; M=12 },
    @{ n=5; A=3; B="no"; C=4; D=3; E=3; F="yes"; G=43711.36485364258; H=2; I=30; J=3; K="no"; L="looking at phone"; M=12 },
    @{ n=6; A=4; B="no"; C=5; D=3; E=3; F="yes"; G=43711.36509579796; H=2; I=30; J=3; K="no"; L="looking at phone"; M=12 },
    @{ n=7; A=5; B="no"; C=6; D=3; E=5; F="yes"; G=43711.36578422158; H=2; I=30; J=2; K="no"; L="meditation";        M=12 },
    @{ n=8; A=6; B="no"; C=7; D=4; E=3; F="yes"; G=43711;              H=2; I=45; J=2; K="no"; L="looking at phone"; M=12 }
)

foreach ($r in $rows) {
    $n = $r.n
    $ws.Cells.Item($n, 1).Value  = $r.A
    $ws.Cells.Item($n, 2).Value  = $r.B
    $ws.Cells.Item($n, 3).Value  = $r.C
    $ws.Cells.Item($n, 4).Value  = $r.D
    $ws.Cells.Item($n, 5).Value  = $r.E
    $ws.Cells.Item($n, 6).Value  = $r.F
    $ws.Cells.Item($n, 7).Value  = $r.G
    $ws.Cells.Item($n, 8).Value  = $r.H
    $ws.Cells.Item($n, 9).Value  = $r.I
    $ws.Cells.Item($n, 10).Value = $r.J
    $ws.Cells.Item($n, 11).Value = $r.K
    $ws.Cells.Item($n, 12).Value = $r.L
    $ws.Cells.Item($n, 13).Value = $r.M
}
